# Fill in the missing "CreatedUser" values (column I) for rows 3-7 on the
# "demoaut" sheet. Each of these cells should mirror the "Username" value
# already present in column F for the same row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("demoaut")

for ($row = 3; $row -le 7; $row++) {
    $ws.Cells.Item($row, 9).Value2 = $ws.Cells.Item($row, 6).Value2
}

# Clear the stale selection on the sheet so it no longer points at I2:I7.
$ws.Range("A1").Select()
